# Applies the cell-value updates for Sheet1 (rows 2-11) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.05
$ws.Range("G2").Value = 4.3
$ws.Range("H2").Value = 2.46
$ws.Range("I2").Value = 2.84
$ws.Range("J2").Value = 2.72
$ws.Range("N2").Value = 2.02
$ws.Range("Q2").Value = 1.69
$ws.Range("S2").Value = 6
$ws.Range("U2").Value = 1.62
$ws.Range("V2").Value = 1.54

# Row 3
$ws.Range("G3").Value = 1.2
$ws.Range("H3").Value = 17
$ws.Range("J3").Value = 9.199999999999999
$ws.Range("K3").Value = 11.5
$ws.Range("N3").Value = 7.6
$ws.Range("W3").Value = 6
$ws.Range("AB3").Value = 18
$ws.Range("AF3").Value = 12.5

# Row 4
$ws.Range("G4").Value = 2.02
$ws.Range("Q4").Value = 2.36
$ws.Range("S4").Value = 4.7
$ws.Range("T4").Value = 1.04
$ws.Range("U4").Value = 1.04
$ws.Range("W4").Value = 1.98
$ws.Range("X4").Value = 11
$ws.Range("Z4").Value = 55
$ws.Range("AB4").Value = 7.6
$ws.Range("AC4").Value = 8.199999999999999
$ws.Range("AF4").Value = 10.5
$ws.Range("AG4").Value = 11.5
$ws.Range("AL4").Value = 1000

# Row 5
$ws.Range("R5").Value = 1.24

# Row 6
$ws.Range("F6").Value = 1.62
$ws.Range("G6").Value = 1.69
$ws.Range("H6").Value = 4.6
$ws.Range("I6").Value = 5.4
$ws.Range("J6").Value = 4.9
$ws.Range("K6").Value = 5.8
$ws.Range("P6").Value = 3.4
$ws.Range("Q6").Value = 1.34
$ws.Range("R6").Value = 1.96
$ws.Range("S6").Value = 1.84
$ws.Range("T6").Value = 1.45
$ws.Range("U6").Value = 2.8
$ws.Range("V6").Value = 1.23
$ws.Range("W6").Value = 2.44
$ws.Range("X6").Value = 46
$ws.Range("Y6").Value = 38
$ws.Range("Z6").Value = 60
$ws.Range("AB6").Value = 19
$ws.Range("AD6").Value = 22
$ws.Range("AE6").Value = 50
$ws.Range("AI6").Value = 50
$ws.Range("AJ6").Value = 20
$ws.Range("AL6").Value = 22
$ws.Range("AM6").Value = 55
$ws.Range("AN6").Value = 5.1

# Row 7
$ws.Range("F7").Value = 4.7
$ws.Range("I7").Value = 1.87
$ws.Range("L7").Value = 1.34
$ws.Range("S7").Value = 3.35
$ws.Range("V7").Value = 2.08
$ws.Range("Y7").Value = 10.5
$ws.Range("Z7").Value = 13.5
$ws.Range("AA7").Value = 24
$ws.Range("AB7").Value = 21
$ws.Range("AC7").Value = 10.5
$ws.Range("AD7").Value = 12.5
$ws.Range("AE7").Value = 24
$ws.Range("AG7").Value = 24
$ws.Range("AH7").Value = 24
$ws.Range("AK7").Value = 85
$ws.Range("AL7").Value = 90
$ws.Range("AO7").Value = 15.5

# Row 8
$ws.Range("G8").Value = 1.75
$ws.Range("I8").Value = 8.4
$ws.Range("J8").Value = 3.35
$ws.Range("L8").Value = 1.01
$ws.Range("V8").Value = 1.14
$ws.Range("W8").Value = 2.32

# Row 9
$ws.Range("O9").Value = 1.37
$ws.Range("P9").Value = 1.87
$ws.Range("Q9").Value = 2.08
$ws.Range("U9").Value = 1.68
$ws.Range("AI9").Value = 190
$ws.Range("AO9").Value = 370

# Row 10
$ws.Range("O10").Value = 1.27
$ws.Range("Q10").Value = 1.85
$ws.Range("S10").Value = 3.1
$ws.Range("AL10").Value = 32

# Row 11
$ws.Range("F11").Value = 2.12
$ws.Range("G11").Value = 2.36
$ws.Range("H11").Value = 3.75
$ws.Range("I11").Value = 4.4
$ws.Range("J11").Value = 2.86
$ws.Range("K11").Value = 3.6
$ws.Range("L11").Value = 1.42
$ws.Range("M11").Value = 1.1
$ws.Range("N11").Value = 2.86
$ws.Range("O11").Value = 1.43
$ws.Range("P11").Value = 1.63
$ws.Range("Q11").Value = 2.28
$ws.Range("R11").Value = 1.23
$ws.Range("S11").Value = 4.4
$ws.Range("T11").Value = 1.01
$ws.Range("V11").Value = 1.29
$ws.Range("W11").Value = 1.74
$ws.Range("X11").Value = 12.5
$ws.Range("Y11").Value = 12.5
$ws.Range("Z11").Value = 29
$ws.Range("AA11").Value = 110
$ws.Range("AB11").Value = 8.199999999999999
$ws.Range("AC11").Value = 7.8
$ws.Range("AD11").Value = 17.5
$ws.Range("AE11").Value = 65
$ws.Range("AF11").Value = 13.5
$ws.Range("AG11").Value = 12
$ws.Range("AH11").Value = 22
$ws.Range("AJ11").Value = 32
$ws.Range("AK11").Value = 29
$ws.Range("AL11").Value = 55
$ws.Range("AN11").Value = 29
$ws.Range("AO11").Value = 95
